# Completed unit testing according to test plan - 7 tests passed.
# Fills in the Unit Test Plan table (rows 7-13: Method Inputs / Condition being
# Tested / Expected Result columns) and the Developer name cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: the exact order these cell values are written in determines the
# order new entries land in the shared-string table, so the sequence below
# is deliberately not simple row-major order - it mirrors how the values
# were actually entered.

# --- Test Case 1 (__init__ : Attributes are set to parameter values.) ---
$ws.Range("E7").Value = "None"
$ws.Range("F7").Value = "account_number=2121 client_number=2222 balance=1000.00 date_created=today() management_fee=2.00"

# --- Test Case 2 (__init__ : management fee has invalid type.) ---
$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = 'account_number=2121 client_number=2222 balance=1000.00 date_created=today() management_fee="not_Mark"'
$ws.Range("G8").Value = 2.5499999999999998
$ws.Range("G8").HorizontalAlignment = -4131

$ws.Range("G7").Value = "Object created"

# --- Test Case 3 (get_service_charges : date created more than 10 years ago) ---
$ws.Range("E9").Value = "management_fee=2.00"
$ws.Range("G9").Value = "BASE_SERVICE_CHARGE"

# --- Test Case 4 (get_service_charges : date created exactly 10 years ago) ---
$ws.Range("F10").Value = "account_number=2121 client_number=2222 balance=1000.00 date_created=(today() - 10 years)"
$ws.Range("E10").Value = "management_fee=2.00"

# --- Test Case 5 (get_service_charges : date created within last 10 years.) ---
$ws.Range("F11").Value = "account_number=2121 client_number=2222 balance=1000.00 date_created=2024-01-01"
$ws.Range("G11").Value = "BASE_SERVICE_CHARGE + 2.00"
$ws.Range("E11").Value = "management_fee=2.00"

$ws.Range("F9").Value = "account_number=2121 client_number=2222 balance=1000.00 date_created=2010-01-01"
$ws.Range("G10").Value = 'BASE_SERVICE_CHARGE and "Management Fee: Waived"'

# --- Test Case 6 (__str__ : displays waived management fee when date created more than 10 years ago.) ---
$ws.Range("E12").Value = "management_fee=2.00"
$ws.Range("F12").Value = "account_number=2121 client_number=2222 balance=1000.00 date_created=2010-01-01"
$ws.Range("G12").Value = "Management Fee: Waived"

# --- Test Case 7 (__str__ : displays  management fee when date created within last 10 years.) ---
$ws.Range("E13").Value = "management_fee=2.00"
$ws.Range("F13").Value = "account_number=2121 client_number=2222 balance=1000.00 date_created=2024-01-01"
$ws.Range("G13").Value = "Management Fee: $2.00"

# Developer name (set last so it lands at the end of the shared-string table)
$ws.Range("C3").Value = "Ralph Vitug"

# Final cursor position left on the sheet when the workbook was saved
$ws.Range("F12").Select() | Out-Null
